$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DESFAZIMENTOS")

# Delete row 2 (the "8.2024.1026/000317-7" / "Aguarda Autorização" /
# "FUNDACAO MARIA TAVARES" / "DG/PRESIDÊNCIA" record); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# Re-apply the AutoFilter so its range shrinks from A1:M8 to A1:M7 (Excel
# does not automatically resize it when rows are deleted via Delete()).
$ws.AutoFilterMode = $false
$ws.Range("A1:M7").AutoFilter()

$ws.Activate()
$ws.Range("C2:C4").Select()
